# ESCALETA_CN_08_01_CO.xlsx
# "Inclusión de proyectos método científico" - adds a new resource row
# (a "Competencias" entry about formulating hypotheses in an experiment)
# right before the existing "Fin de unidad" row on the ESCALETA sheet.
#
# Net effect on the ESCALETA sheet: a new row is inserted at row 30,
# pushing the former rows 30-32 down to 31-33, and the (now vacated)
# row 30 is populated with the new "Competencias" resource data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESCALETA")

# --- Shift the existing rows 30-32 down to 31-33 -----------------------
# Range.Copy() preserves the destination's existing style ids (which, on
# this sheet, are already identical row-by-row), but it does NOT clear
# cells that are non-blank in the destination while blank in the source,
# so each destination row is cleared first.

$ws.Range("A33:U33").ClearContents()
$ws.Range("A32:U32").Copy($ws.Range("A33:U33"))

$ws.Range("A32:U32").ClearContents()
$ws.Range("A31:U31").Copy($ws.Range("A32:U32"))

$ws.Range("A31:U31").ClearContents()
$ws.Range("A30:U30").Copy($ws.Range("A31:U31"))

# "Número" column is sequential by row position, independent of content.
$ws.Range("H30").Value = 28
$ws.Range("H31").Value = 29
$ws.Range("H32").Value = 30
$ws.Range("H33").Value = 31

# --- Populate row 30 with the new "Competencias" resource ---------------
# Columns A, B, C, E, F, H, I, K, N, O keep the values already in row 30
# (Asignatura, Código del guion, Título del guion, Número, Fichas,
# Tipología, Comentarios/Indicaciones, Aparece en Cuaderno all stay the
# same for this guion).
$ws.Range("D30").Value = "Competencias"
$ws.Range("G30").Value = "Competencias: plantear hipótesis en un experimento"
$ws.Range("J30").Value = "Proyecto para generar preguntas e hipótesis en un experimento científico"
$ws.Range("L30").Value = "INTERACTIVO "
$ws.Range("M30").Value = "F13"
$ws.Range("P30").Value = "NO"
$ws.Range("Q30").Value = 6
$ws.Range("R30").Value = "RF"
$ws.Range("S30").Value = "Recursos F"
$ws.Range("T30").Value = "Recurso F13-02"
$ws.Range("U30").Value = "RF_01_01_CO"
